# Fix libre file parsing: the Boolean columns in all three sheets were
# accidentally carrying a date-style cell format (style index 1 -> numFmtId 14)
# inherited from the preceding DateTime column. Re-apply the correct
# "General" number format to those Boolean cells so they get their own
# style (numFmtId 0 / General, applyNumberFormat=1) instead of a date format.

$wb = $excel.ActiveWorkbook

$wsWithTable = $wb.Worksheets.Item("WithTable")
$wsWithTable.Range("D2:D5").NumberFormat = "General"

$wsTableless = $wb.Worksheets.Item("Tableless")
$wsTableless.Range("D2:D5").NumberFormat = "General"

$wsDuplicate = $wb.Worksheets.Item("WithTable_Duplicate")
$wsDuplicate.Range("E5:E8").NumberFormat = "General"

# Update each sheet's remembered selection.
$null = $wsWithTable.Range("E23").Select()
$null = $wsTableless.Range("I19").Select()
$null = $wsDuplicate.Range("M11").Select()

# Move the active/visible tab from "Tableless" to "WithTable_Duplicate".
$null = $wsDuplicate.Activate()
$null = $wsDuplicate.Range("M11").Select()
